# Scheduled-runner style refresh of market-price-derived columns
# (currentAveragePrice*, Leve*Price*, Leve*Profit*) across the per-job
# leve worksheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR). Only data values are
# touched - no formulas, rows, or formatting are involved.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 3018.0435
$ws.Range("I5").Value = 814.75
$ws.Range("J5").Value = 8054.143
$ws.Range("K5").Value = 814.75
$ws.Range("L5").Value = 8054.143
$ws.Range("M5").Value = -699.75
$ws.Range("N5").Value = -8284.143
$ws.Range("H40").Value = 7214.778
$ws.Range("J40").Value = 7247.25
$ws.Range("L40").Value = 7247.25
$ws.Range("N40").Value = -7597.25
$ws.Range("H74").Value = 4931500
$ws.Range("I74").Value = 11908668
$ws.Range("J74").Value = 6440.5293
$ws.Range("K74").Value = 11908668
$ws.Range("L74").Value = 6440.5293
$ws.Range("M74").Value = -11907732
$ws.Range("N74").Value = -8312.5293
$ws.Range("H77").Value = 4931500
$ws.Range("I77").Value = 11908668
$ws.Range("J77").Value = 6440.5293
$ws.Range("K77").Value = 59543340
$ws.Range("L77").Value = 32202.6465
$ws.Range("M77").Value = -59538660
$ws.Range("N77").Value = -41562.6465
$ws.Range("H86").Value = 333334820
$ws.Range("I86").Value = 333334820
$ws.Range("K86").Value = 333334820
$ws.Range("M86").Value = -333333697
$ws.Range("H89").Value = 333334820
$ws.Range("I89").Value = 333334820
$ws.Range("K89").Value = 1666674100
$ws.Range("M89").Value = -1666668484
$ws.Range("H113").Value = 11684.692
$ws.Range("I113").Value = 5989.8335
$ws.Range("J113").Value = 16566
$ws.Range("K113").Value = 5989.8335
$ws.Range("L113").Value = 16566
$ws.Range("M113").Value = -2735.8335
$ws.Range("N113").Value = -23074
$ws.Range("H138").Value = 4710.793
$ws.Range("J138").Value = 5409.875
$ws.Range("L138").Value = 16229.625
$ws.Range("N138").Value = -26509.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5302.032
$ws.Range("J2").Value = 12205.2
$ws.Range("L2").Value = 12205.2
$ws.Range("N2").Value = -12431.2
$ws.Range("H61").Value = 2982.9524
$ws.Range("I61").Value = 1935.4445
$ws.Range("J61").Value = 3768.5833
$ws.Range("K61").Value = 1935.4445
$ws.Range("L61").Value = 3768.5833
$ws.Range("M61").Value = -1723.4445
$ws.Range("N61").Value = -4192.5833
$ws.Range("H116").Value = 5302.032
$ws.Range("J116").Value = 12205.2
$ws.Range("L116").Value = 12205.2
$ws.Range("N116").Value = -16793.2
$ws.Range("H136").Value = 2982.9524
$ws.Range("I136").Value = 1935.4445
$ws.Range("J136").Value = 3768.5833
$ws.Range("K136").Value = 5806.333500000001
$ws.Range("L136").Value = 11305.7499
$ws.Range("M136").Value = -3256.333500000001
$ws.Range("N136").Value = -16405.7499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5302.032
$ws.Range("J3").Value = 12205.2
$ws.Range("L3").Value = 12205.2
$ws.Range("N3").Value = -12433.2
$ws.Range("H22").Value = 2501573
$ws.Range("I22").Value = 2150
$ws.Range("K22").Value = 2150
$ws.Range("M22").Value = -1977
$ws.Range("H100").Value = 33099.832
$ws.Range("J100").Value = 33099.832
$ws.Range("L100").Value = 33099.832
$ws.Range("N100").Value = -35263.832
$ws.Range("H105").Value = 3750.923
$ws.Range("J105").Value = 1331.6666
$ws.Range("L105").Value = 1331.6666
$ws.Range("N105").Value = -4825.6666
$ws.Range("H107").Value = 11112944
$ws.Range("I107").Value = 16668079
$ws.Range("K107").Value = 16668079
$ws.Range("M107").Value = -16666159

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 31255858
$ws.Range("I134").Value = 111112856
$ws.Range("J134").Value = 7465.522
$ws.Range("K134").Value = 333338568
$ws.Range("L134").Value = 22396.566
$ws.Range("M134").Value = -333336033
$ws.Range("N134").Value = -27466.566

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4166.1665
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 4599.4
$ws.Range("K80").Value = 6000
$ws.Range("L80").Value = 13798.2
$ws.Range("M80").Value = -5064
$ws.Range("N80").Value = -15670.2
$ws.Range("H83").Value = 4166.1665
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 4599.4
$ws.Range("K83").Value = 18000
$ws.Range("L83").Value = 41394.6
$ws.Range("M83").Value = -13320
$ws.Range("N83").Value = -50754.6
$ws.Range("H113").Value = 417453.9
$ws.Range("J113").Value = 1001196
$ws.Range("L113").Value = 3003588
$ws.Range("N113").Value = -3007928
$ws.Range("H131").Value = 36113332
$ws.Range("J131").Value = 19610614
$ws.Range("L131").Value = 58831842
$ws.Range("N131").Value = -58841922
$ws.Range("H134").Value = 4997.4814
$ws.Range("I134").Value = 4997.4814
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 14992.4442
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -9922.444199999998
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2582.15
$ws.Range("I80").Value = 1224
$ws.Range("J80").Value = 8014.75
$ws.Range("K80").Value = 1224
$ws.Range("L80").Value = 8014.75
$ws.Range("M80").Value = -226
$ws.Range("N80").Value = -10010.75
$ws.Range("H83").Value = 2582.15
$ws.Range("I83").Value = 1224
$ws.Range("J83").Value = 8014.75
$ws.Range("K83").Value = 6120
$ws.Range("L83").Value = 40073.75
$ws.Range("M83").Value = -1128
$ws.Range("N83").Value = -50057.75
$ws.Range("H122").Value = 7682.8667
$ws.Range("I122").Value = 6955.5713
$ws.Range("J122").Value = 8319.25
$ws.Range("K122").Value = 20866.7139
$ws.Range("L122").Value = 24957.75
$ws.Range("M122").Value = -18416.7139
$ws.Range("N122").Value = -29857.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5929.1875
$ws.Range("I16").Value = 1986.8
$ws.Range("J16").Value = 12499.833
$ws.Range("K16").Value = 1986.8
$ws.Range("L16").Value = 12499.833
$ws.Range("M16").Value = -1816.8
$ws.Range("N16").Value = -12839.833
$ws.Range("H22").Value = 524.625
$ws.Range("J22").Value = 499.5
$ws.Range("L22").Value = 499.5
$ws.Range("N22").Value = -1089.5
$ws.Range("H27").Value = 524.625
$ws.Range("J27").Value = 499.5
$ws.Range("L27").Value = 499.5
$ws.Range("N27").Value = -713.5
$ws.Range("H55").Value = 3001
$ws.Range("I55").Value = 790.75
$ws.Range("J55").Value = 4965.6665
$ws.Range("K55").Value = 790.75
$ws.Range("L55").Value = 4965.6665
$ws.Range("M55").Value = -617.75
$ws.Range("N55").Value = -5311.6665
$ws.Range("H61").Value = 5849.769
$ws.Range("I61").Value = 3849.7
$ws.Range("K61").Value = 3849.7
$ws.Range("M61").Value = -3647.7
$ws.Range("H68").Value = 2181.4546
$ws.Range("I68").Value = 2181.4546
$ws.Range("K68").Value = 2181.4546
$ws.Range("M68").Value = -1432.4546
$ws.Range("H71").Value = 2181.4546
$ws.Range("I71").Value = 2181.4546
$ws.Range("K71").Value = 10907.273
$ws.Range("M71").Value = -7163.273000000001
$ws.Range("H113").Value = 5849.769
$ws.Range("I113").Value = 3849.7
$ws.Range("K113").Value = 3849.7
$ws.Range("M113").Value = -1679.7
$ws.Range("H136").Value = 16668936
$ws.Range("I136").Value = 33334262
$ws.Range("J136").Value = 3610.3333
$ws.Range("K136").Value = 100002786
$ws.Range("L136").Value = 10830.9999
$ws.Range("M136").Value = -100000236
$ws.Range("N136").Value = -15930.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 7689.5
$ws.Range("I126").Value = 4004
$ws.Range("K126").Value = 12012
$ws.Range("M126").Value = -9542
